# Singapore GP comparison table: update simulation-vs-actual error metrics
# and re-sync the George Russell / Charles Leclerc row ordering with the
# refreshed notebook output (row/col indices are 1-based; row 1 is the
# header row "driver_name | laps_completed | position_sim | ... | gap_error").
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 2: Lando Norris
$t.Cell(2,5).Range.Text = "-5"
# Row 3: Max Verstappen
$t.Cell(3,3).Range.Text = "-2"
$t.Cell(3,5).Range.Text = "-1"

# Row 4: Oscar Piastri
$t.Cell(4,6).Range.Text = "238.906"
$t.Cell(4,7).Range.Text = "10.221"

# Row 5 used to be Charles Leclerc, now holds George Russell's data
# (the two drivers swapped table positions).
$t.Cell(5,1).Range.Text = "George Russell"
$t.Cell(5,3).Range.Text = "5"
$t.Cell(5,4).Range.Text = "4"
$t.Cell(5,5).Range.Text = "0"
$t.Cell(5,6).Range.Text = "230.048"
$t.Cell(5,7).Range.Text = "1.362"

# Row 6 used to be George Russell, now holds Charles Leclerc's data.
$t.Cell(6,1).Range.Text = "Charles Leclerc"
$t.Cell(6,3).Range.Text = "4"
$t.Cell(6,4).Range.Text = "5"
$t.Cell(6,5).Range.Text = "-2"
$t.Cell(6,6).Range.Text = "236.344"
$t.Cell(6,7).Range.Text = "7.659"

# Row 7: Lewis Hamilton
$t.Cell(7,5).Range.Text = "7"
$t.Cell(7,6).Range.Text = "215.298"
$t.Cell(7,7).Range.Text = "-13.388"

# Row 8: Carlos Sainz
$t.Cell(8,5).Range.Text = "1"
$t.Cell(8,6).Range.Text = "220.531"
$t.Cell(8,7).Range.Text = "-8.154"

# Row 9: Fernando Alonso
$t.Cell(9,6).Range.Text = "327.238"
$t.Cell(9,7).Range.Text = "0.631"

# Row 10: Franco Colapinto
$t.Cell(10,6).Range.Text = "326.161"
$t.Cell(10,7).Range.Text = "-0.446"

# Row 11: Nico Hulkenberg
$t.Cell(11,6).Range.Text = "328.364"
$t.Cell(11,7).Range.Text = "1.756"

# Row 12: Sergio Perez
$t.Cell(12,5).Range.Text = "-2"

# Row 13: Yuki Tsunoda
$t.Cell(13,6).Range.Text = "330.731"
$t.Cell(13,7).Range.Text = "4.123"

# Row 14: Lance Stroll
$t.Cell(14,6).Range.Text = "318.564"
$t.Cell(14,7).Range.Text = "-8.044"

# Row 15: Esteban Ocon
$t.Cell(15,6).Range.Text = "324.342"
$t.Cell(15,7).Range.Text = "-2.266"

# Row 16: Pierre Gasly
$t.Cell(16,5).Range.Text = "-6"
$t.Cell(16,6).Range.Text = "315.917"
$t.Cell(16,7).Range.Text = "-10.690"

# Row 17: Guanyu Zhou
$t.Cell(17,6).Range.Text = "326.337"
$t.Cell(17,7).Range.Text = "-0.271"

# Row 18: Valtteri Bottas
$t.Cell(18,6).Range.Text = "327.224"
$t.Cell(18,7).Range.Text = "0.616"

# Row 19: Daniel Ricciardo
$t.Cell(19,6).Range.Text = "361.258"
$t.Cell(19,7).Range.Text = "34.650"

# Row 20: Alexander Albon (overtake_error was "nan")
$t.Cell(20,5).Range.Text = "-3"

# Row 21: Kevin Magnussen (overtake_error was "nan")
$t.Cell(21,5).Range.Text = "-4"
